# Add results row for network/dijkstra benchmark
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "network/dijkstra/runme_large.sh"
$ws.Range("B10").Value = 0.02
$ws.Range("C10").Value = 0.01
$ws.Range("D10").Value = 0

# Match the author's final cursor position recorded in the saved file
[void]$ws.Range("A27").Select()
